$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.518.81"
$ws.Range("E2").Value = "  -3.83%  "

$ws.Range("D3").Value = "3.191.36"
$ws.Range("E3").Value = "  -4.84%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'533.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.75%  "

$ws.Range("D6").Value = "'134.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.38%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.189.58"

$ws.Range("E9").Value = "  -4.90%  "

$ws.Range("D10").Value = "'7.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.49%  "

$ws.Range("E11").Value = "  -7.07%  "

$ws.Range("E12").Value = "  -5.02%  "

$ws.Range("D13").Value = "3.738.49"
$ws.Range("E13").Value = "  -4.93%  "

$ws.Range("E14").Value = "  -0.46%  "

$ws.Range("D15").Value = "'25.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.88%  "

$ws.Range("D16").Value = "3.195.35"
$ws.Range("E16").Value = "  -4.95%  "

$ws.Range("D17").Value = "58.577.03"
$ws.Range("E17").Value = "  -3.85%  "

$ws.Range("E18").Value = "  -7.63%  "

$ws.Range("D20").Value = "'13.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.38%  "

$ws.Range("E21").Value = "  -8.69%  "

$ws.Range("D22").Value = "'358.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.68%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E24").Value = "  -7.04%  "

$ws.Range("D25").Value = "'0.517"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.79%  "

$ws.Range("E26").Value = "  -5.05%  "

$ws.Range("D27").Value = "'0.170"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.15%  "

$ws.Range("D28").Value = "0.0₃0952"
$ws.Range("E28").Value = "  -11.80%  "

$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").Value = "'7.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.81%  "

$ws.Range("E32").Value = "  -8.11%  "

$ws.Range("D33").Value = "'7.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.89%  "

$ws.Range("D34").Value = "'21.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.11%  "

$ws.Range("E35").Value = "  -9.09%  "

$ws.Range("E36").Value = "  -7.25%  "

$ws.Range("E37").Value = "  -4.72%  "

$ws.Range("E38").Value = "  -7.14%  "

$ws.Range("E39").Value = "  -8.19%  "

$ws.Range("D40").Value = "'25.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.00%  "

$ws.Range("D41").Value = "'0.0705"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.04%  "

$ws.Range("D42").Value = "3.221.01"
$ws.Range("E42").Value = "  -5.03%  "

$ws.Range("D43").Value = "'40.74"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'0.708"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.72%  "

$ws.Range("E45").Value = "  -3.75%  "

$ws.Range("E46").Value = "  -6.65%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.51%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("D49").Value = "2.280.08"
$ws.Range("E49").Value = "  -8.61%  "

$ws.Range("D50").Value = "'6.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.06%  "

$ws.Range("D51").Value = "'20.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.43%  "
